# Cross-sheet fixture update: "Actually evaluate open ended ranges."
#
# Sheet1 gains a header row (shifting existing data down by one row),
# plus a new column D whose second row sums the now open-ended column A
# range (SUM(Sheet1!A:A)). Sheet2's cross-sheet formulas, and the
# workbook-level defined name, follow the row shift automatically /
# are updated to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# 1. Insert a new row above row 1 on Sheet1; existing formulas/refs
#    (same-sheet and cross-sheet) shift automatically.
$ws1.Rows.Item(1).Insert()

# 2. Populate the new header row with labels.
$ws1.Range("A1").Value = "Hello"
$ws1.Range("B1").Value = "This "
$ws1.Range("C1").Value = "Is"
$ws1.Range("D1").Value = "A Table"

# 3. Add the open-ended SUM formula in the new column D, row 2.
$ws1.Range("D2").Formula = "=SUM(Sheet1!A:A)"

# 4. Update the workbook-level defined name to track the shifted cell.
$wb.Names.Item("LastCell").RefersTo = "=Sheet1!`$C`$5"

# 5. Restore view/selection state: Sheet1 becomes the active sheet with
#    selection at M36; Sheet2 keeps its D1 selection but is no longer
#    the active tab.
$ws1.Activate()
$ws1.Range("M36").Select()
$ws2.Range("D1").Select()
$ws1.Activate()

# 6. Column-metrics nudge recorded alongside this edit.
$ws1.StandardWidth = 8.54296875
$ws2.StandardWidth = 8.54296875

# 7. Re-stamp the pre-existing data cells with an explicit "Normal"
#    style so they carry their own style record, matching the
#    original cells' appearance (the new header/SUM cells keep the
#    sheet's default/unstyled look).
$ws1.Range("A1:C5").Style = "Normal"
$ws2.Range("A1:C3").Style = "Normal"
